$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3430825769901276
$ws.Range("B1").Value = 3.416850566864014
$ws.Range("C1").Value = 5.935484409332275
$ws.Range("D1").Value = 1.657985687255859
$ws.Range("E1").Value = 0.9946198463439941
